$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 52 (start of Maze Size 15 block)
$ws.Cells.Item(52, 1).Value = 15
$ws.Cells.Item(52, 2).Value = 107
$ws.Cells.Item(52, 3).Formula = "=AVERAGE(B52:B60)"
$ws.Cells.Item(52, 4).Value = 0.1206

$ws.Cells.Item(53, 2).Value = 119
$ws.Cells.Item(54, 2).Value = 135
$ws.Cells.Item(55, 2).Value = 116
$ws.Cells.Item(56, 2).Value = 126
$ws.Cells.Item(57, 2).Value = 130
$ws.Cells.Item(58, 2).Value = 122
$ws.Cells.Item(59, 2).Value = 129
$ws.Cells.Item(60, 2).Value = 101

# Row 61 (start of Maze Size 20 block)
$ws.Cells.Item(61, 1).Value = 20
$ws.Cells.Item(61, 2).Value = 108
$ws.Cells.Item(61, 3).Formula = "=AVERAGE(B61:B70)"
$ws.Cells.Item(61, 4).Value = 0.1136

$ws.Cells.Item(62, 2).Value = 117
$ws.Cells.Item(63, 2).Value = 127
$ws.Cells.Item(64, 2).Value = 126
$ws.Cells.Item(65, 2).Value = 101
$ws.Cells.Item(66, 2).Value = 101
$ws.Cells.Item(67, 2).Value = 123
$ws.Cells.Item(68, 2).Value = 119
$ws.Cells.Item(69, 2).Value = 110
$ws.Cells.Item(70, 2).Value = 104

# Update view: scroll so A46 is the top-left visible cell, and select E66
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 46
$win.ScrollColumn = 1
$ws.Range("E66").Select()
